# Add a new CV entry row for "COE Quantitative Curriculum Review Committee
# Member" service, inserted as the new row 18 (right after the existing
# "Faculty Advisory Committee..." service row), pushing all subsequent
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 18; formatting (styles for C/D columns) is
# inherited from the row above, matching the rest of the "service" block.
$ws.Range("A18").EntireRow.Insert()

# Match the row height used by the sibling "service" rows (16/17/19).
$ws.Rows.Item(18).RowHeight = 34

# Populate the new service entry.
$ws.Range("A18").Value = "service"
$ws.Range("B18").Value = 2019
$ws.Range("C18").Value = "current"
$ws.Range("D18").Value = "COE Quantitative Curriculum Review Committee Member"
$ws.Range("E18").Value = "University of Oregon"

# Leave the selection where the author's last edit left it.
$ws.Range("G18").Select() | Out-Null
